# Settings back to qfin22 presentation - recalibration
$wb = $excel.ActiveWorkbook

# --- Linear sheet (mu / B / sig2 / abs_epsi_autocorr) ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.1458682256149341
$wsLinear.Range("B3").Value = 0.004274657799904361
$wsLinear.Range("B4").Value = 591.6082330581754
$wsLinear.Range("B5").Value = "[1.0, 0.14972786055093568, -0.04415492273241662, 0.0009662612215151185, -0.030944553166090707, -0.040943670894978954, 0.1412589302918169, 0.28528686545182336, 0.14563304500915772, -0.0631672013487196, -0.024044140765320042, -0.03119380771567881, -0.056257310939128176, 0.12848088964105359, 0.28585663365154634, 0.09528724559701804, -0.05358760087002149, -0.04687610456612523, -0.04047555426666085, -0.06860976474910684]"

# --- NonLinear sheet (mu_0 / B_0 / sig2_0 / mu_1 / B_1 / sig2_1 / abs_epsi_autocorr) ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -1.300070168473521
$wsNonLinear.Range("B5").Value = -0.06909160758583491
$wsNonLinear.Range("B6").Value = 581.1261898492172
$wsNonLinear.Range("B7").Value = 0.9033016195215622
$wsNonLinear.Range("B8").Value = -0.06715249279225555
$wsNonLinear.Range("B9").Value = 600.8476745695534
$wsNonLinear.Range("B10").Value = "[1.0, 0.14613125777637287, -0.039486963568784526, 0.004521675873399327, -0.02707178504151291, -0.03788434436203105, 0.14017510079564816, 0.2778376587034649, 0.14437620219647593, -0.0584675595365691, -0.01998631367510322, -0.026930653836475433, -0.051487484108415846, 0.1272245635457431, 0.2790305678553474, 0.09238020012227396, -0.04930088178482296, -0.0448933825800163, -0.03802305876274631, -0.06413476336464809]"
